$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-24 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-25 Monday", 2) | Out-Null
$d.Content.Find.Execute("30×40=", $true, $false, $false, $false, $false, $true, 1, $false, "49×14=", 2) | Out-Null
$d.Content.Find.Execute("18×57=", $true, $false, $false, $false, $false, $true, 1, $false, "49×36=", 2) | Out-Null
$d.Content.Find.Execute("56×19=", $true, $false, $false, $false, $false, $true, 1, $false, "29×12=", 2) | Out-Null
$d.Content.Find.Execute("36×50=", $true, $false, $false, $false, $false, $true, 1, $false, "89×77=", 2) | Out-Null
$d.Content.Find.Execute("44×87=", $true, $false, $false, $false, $false, $true, 1, $false, "75×42=", 2) | Out-Null
$d.Content.Find.Execute("52×43=", $true, $false, $false, $false, $false, $true, 1, $false, "59×53=", 2) | Out-Null
$d.Content.Find.Execute("53×31=", $true, $false, $false, $false, $false, $true, 1, $false, "54×92=", 2) | Out-Null
$d.Content.Find.Execute("98×59=", $true, $false, $false, $false, $false, $true, 1, $false, "16×38=", 2) | Out-Null
$d.Content.Find.Execute("31×89=", $true, $false, $false, $false, $false, $true, 1, $false, "33×94=", 2) | Out-Null
$d.Content.Find.Execute("45×74=", $true, $false, $false, $false, $false, $true, 1, $false, "29×46=", 2) | Out-Null
$d.Content.Find.Execute("36×26=", $true, $false, $false, $false, $false, $true, 1, $false, "15×27=", 2) | Out-Null
$d.Content.Find.Execute("12×53=", $true, $false, $false, $false, $false, $true, 1, $false, "12×13=", 2) | Out-Null
$d.Content.Find.Execute("53×13=", $true, $false, $false, $false, $false, $true, 1, $false, "36×39=", 2) | Out-Null
$d.Content.Find.Execute("46×62=", $true, $false, $false, $false, $false, $true, 1, $false, "25×68=", 2) | Out-Null
$d.Content.Find.Execute("51×40=", $true, $false, $false, $false, $false, $true, 1, $false, "27×73=", 2) | Out-Null
$d.Content.Find.Execute("78×81=", $true, $false, $false, $false, $false, $true, 1, $false, "45×91=", 2) | Out-Null
$d.Content.Find.Execute("26×52=", $true, $false, $false, $false, $false, $true, 1, $false, "36×85=", 2) | Out-Null
$d.Content.Find.Execute("56×56=", $true, $false, $false, $false, $false, $true, 1, $false, "67×25=", 2) | Out-Null
$d.Content.Find.Execute("16×25=", $true, $false, $false, $false, $false, $true, 1, $false, "79×71=", 2) | Out-Null
$d.Content.Find.Execute("42×99=", $true, $false, $false, $false, $false, $true, 1, $false, "32×45=", 2) | Out-Null
$d.Content.Find.Execute("77×91=", $true, $false, $false, $false, $false, $true, 1, $false, "52×73=", 2) | Out-Null
$d.Content.Find.Execute("93×46=", $true, $false, $false, $false, $false, $true, 1, $false, "96×78=", 2) | Out-Null
$d.Content.Find.Execute("64×43=", $true, $false, $false, $false, $false, $true, 1, $false, "15×53=", 2) | Out-Null
$d.Content.Find.Execute("66×45=", $true, $false, $false, $false, $false, $true, 1, $false, "64×30=", 2) | Out-Null
$d.Content.Find.Execute("37×56=", $true, $false, $false, $false, $false, $true, 1, $false, "52×61=", 2) | Out-Null
